# Fruta / hortaliza, semanal
#
# A new weekly price-report row is inserted at row 357 of the sheet
# (pushing the existing rows 357-421 down to 358-422), and populated with
# the new record's values. This mirrors Excel's own "insert row, shift
# cells down" behaviour, which is what the canonical OOXML diff shows
# (dimension grows from A1:R421 to A1:R422, and every former row N in
# 357..421 becomes row N+1 with identical contents).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 357, shifting 357:421 -> 358:422.
$ws.Rows("357:357").Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A357").Value = 9
$ws.Range("B357").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C357").Value = "Metropolitana"
$ws.Range("D357").Value = 45244
$ws.Range("E357").Value = 13
$ws.Range("F357").Value = 100112026
$ws.Range("G357").Value = "Haba"
$ws.Range("H357").Value = "Sin especificar"
$ws.Range("I357").Value = "Primera"
$ws.Range("J357").Value = 160
$ws.Range("K357").Value = 9000
$ws.Range("L357").Value = 10000
$ws.Range("M357").Value = 9500
$ws.Range("N357").Value = "$/saco 25 kilos"
$ws.Range("O357").Value = "Región Metropolitana"
$ws.Range("P357").Value = 380
$ws.Range("Q357").Value = 25
$ws.Range("R357").Value = "Hortaliza"
